$d = $word.ActiveDocument
$replacements = @(
    ,@("Video Title", "የቪዲዮ አርዕስት")
    ,@("Ant's Proble", "የጉንዳን ጥያቄ")
    ,@("Topic", "የትምህርት ርዕስ")
    ,@("Problem solving", "ችግር መፍታት")
    ,@("Aim(s)", "አላማ(ዎች)")
    ,@("Finding out how a change in perspective can turn a difficult problem into an easy one.", "የአመለካከታችን መቀየር አስቸጋር ጥያቄዎችን እነዴት ወደቀላል ሊቀይረው እንደሚችል መሳየት፡፡")
    ,@("Length", "ርዝመት")
    ,@("Camp Location", "የካምፕ አድራሻ")
    ,@("Facilitators", "አስተባባሪዎች")
    ,@("N. of students", "የ ተማሪዎች ብዛት")
    ,@("Date", "ቀን")
    ,@("Resources", "አስፈላጊ መሳሪያዎች")
    ,@("needed", "አስፈላጊ")
    ,@("Pen and Paper", "እስኪርብቶ እና ወረቀት")
    ,@("Preparations", "ዠግጁነቶች")
    ,@("None", "ምንም")
    ,@("Video time", "የቪዲዮ ጊዜ")
    ,@("What facilitator does", "የአስተባባሪ ስራ")
    ,@("What learners do", "የተማሪ ስራ")
    ,@("General VMC Video Introduction", "አጠቃላይ የቨርቹዋል ሒሳባ ካምፕ ቪዲዮ ማሰታወቂያ")
    ,@("Video Introduction", "ቪዲዮ ማሰታወቂያ")
    ,@("Riddle", "እንቆቅልሽ")
    ,@("Assist the process, provoke thoughts", "ሂደቱን ማገዝ፣ ሃሳብን ማነሳሳት")
    ,@("3 Ants version", "የ 3 ጉንዳን ስሪት")
    ,@("Solution", "መፍትሄ")
    ,@("This problem is about a change in perspective:", "ይህ ጥያቄ አመለካከት ስለመቀየር ነው፡-")
    ,@("Imagine the same problem but with a single difference in the statement: the ants do not bounce and change verse when they collide, but rather walk on top of each other and keep on moving as if nothing happened.", "አስቡ አንድ አይነት ጥያቄ ግን በትነሽ የአገላለጽ ልዩንት፡-ጉነዳኖቹ  ሲጋጩ አይነሱም አካሄድ አይቀይሩም፣ ግን አንዱ ባነዱ ላይ ይሄዳል በዚሁ የቀጠላሉ ምንም ካልተፈጠረ፡፡")
    ,@("If you think about this second statement, you will notice that the problem is not really changing:", "ስለሁተኛው አረፍተነገር ብታስቡ፣ ጥቄው እነዳልተቀየረ ትረዳላችሁ፡-")
    ,@("If you watch points move on a segment, it is impossible to distinguish between bouncing points and surpassing points.", "በቁራጩ ላይ የሚሄዱትን ነጠቦች በታዩ፣ የመነሳት ነጥቦችንና በላይ የማለፍ ነጥቦችን መለየት ይቻላል፡፡")
    ,@("Are you able to tell if the image above was created thinking about bouncing ants or surpassing ants?", "ከላይ ያለው ምስል ሃሳብ ከፈጠረባችሁ ስለ ጉንዳኖች መነሳት ወይም የጉንዳኖች  በላይ ማለፍ መናገር ተችላላችሁ? ")
    ,@("If you look closely you will realize that it is impossible to tell.", "በቅርበት ካያችሁት እነደሚቻል ትረዳላችሁ፡፡")
    ,@("This:", "ይህ፡-")
    ,@("and this:", "እና ይህ")
    ,@("Are indistinguishable unless you name the ants ( A,B and C in the example)", "የሚለዩ አይደሉም ጉነዳኖቹን ስም ካልሰጣችኋቸው በስተቀር ( A,B እና C በ ምሳሌው)")
    ,@("BUT", "ነገር ግን")
    ,@("By just having a different way of stating the same problem, finding the solution is now easy:", "ለአንድ አይነት ጥያቆ የተለያየ አገላለጽ መኖር መፍትሄውን ለማገኘት ቀላል ያረገዋል፡-")
    ,@("Each ant will walk straight until it falls from one edge. Meaning that every initial position of the ants (no matter how many ants) will last the most if one ant starts from an edge walking towards the other edge.", "እያነዳነዱ ጉንዳን በቀጥታ መስመር የጓዛል ከአንዱ ጫፍ እስከሚዎድቅ ድረስ፡፡ ማለትም እያንዳንዱ የጉንዳኑ መነሻ (ምንም ያክል ጉንዳኖች ቢኖሩ) አንድ ጉንዳን ከአንድ ጫፍ ተነስ እስከ ሌለኛው ጫፍ ከተጓዘ ፈጥኖ ያበቃል፡፡")
    ,@("If you think about surpassing ants this means that the ant that starts the furthest away from the edge, is the last one to fall.", "ስለ ጉንዳኖች ከላይ ማለፍ ካሰብን ይህ ማለት ከሩቅ ጫፍ የጀመረው ጉንዳን መጨረሻ የሚወድቀው ነው፡፡")
    ,@("If you think about bouncing ants, you still don’t know which specific ant will be the last to fall, but if an ant starts from the edge you know that there exists an ant that will fall after 1 meter of walk.", "ስለ ጉንዳኖች መነሳት ካሰብን፣ የመጨረሰ,ሻው የሚወድቀው ጉንዳን የትኛው እነደሆነ አሁንም ማወቅ አንችልም፣ ግን የሆነ ጉንዳን ካነዱ ጫፍ ከጀመረ ከ እንድ ሜትር ጉዞ በኋላ የሚዎድቅ ጉንዳን እነዳለ ታወቃላችሁ፡፡")
    ,@("So, it doesn’t matter how many ants are on the cliff or how they are positioned at the start, as long as one ant is starting from the edge facing the other edge.", "ስለሆነም፣ ስንት ጉንዳኖች በ ጠርዙ ላይ እነዳሉ ወይም ሲጀመር የት እነደተቀመጡ ለውጥ አያመጣም፣ አንድ ጉንዳን ከአንድ ጫፍ ተነስቶ ወደሌለኛው ጫፍ ከተጓዘ፡፡")
)

$count = 0
foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if ($found) { $count++ }
    else { Write-Output "NOT FOUND: $old" }
}
Write-Output "Replaced $count of $($replacements.Count)"
